$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 changes from a "Mean" row of AVERAGE(...) formulas into a "Count" row
# of plain static values.
$ws.Range("A13").Value = "Count"
$ws.Range("B13").Value = 11
$ws.Range("C13").Value = 11
$ws.Range("D13").Value = 3
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 7
$ws.Range("G13").Value = 9
$ws.Range("H13").Value = 9
$ws.Range("I13").Value = 8
$ws.Range("J13").Value = 8
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 5
$ws.Range("M13").Value = 7
$ws.Range("N13").Value = 8
$ws.Range("O13").Value = 6

# Row 14 ("Standard Deviation" / STDEV.S(...) formulas) is wiped out, leaving
# only the blank, normal-height spacer row that row 14 used to be (formerly
# row 15).
$ws.Range("B14:O14").Clear()
$ws.Range("A14").ClearContents()
$ws.Rows.Item(14).AutoFit()

# The old blank spacer row (row 15) is removed entirely so the correlation
# legend rows shift up by one.
$ws.Rows.Item(15).Delete()

# The active selection moves from L13 to A13.
$ws.Range("A13").Select()
